# Weekly price-data update: a new daily record (row 94) is inserted into the
# "Berenjena" sheet, pushing every subsequent record down by one row
# (old row 94 -> 95, ..., old row 131 -> 132).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 94, shifting rows 94:131 -> 95:132.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new observation.
$ws.Range("A94").Value = 6
$ws.Range("B94").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C94").Value = "Metropolitana"
$ws.Range("D94").Value = 44523
$ws.Range("E94").Value = 13
$ws.Range("F94").Value = 100112001
$ws.Range("G94").Value = "Berenjena"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 180
$ws.Range("K94").Value = 10000
$ws.Range("L94").Value = 12000
$ws.Range("M94").Value = 11111
$ws.Range("N94").Value = "$/caja 60 unidades"
$ws.Range("O94").Value = "Provincia de Huasco"
$ws.Range("P94").Value = 185
$ws.Range("Q94").Value = 60
$ws.Range("R94").Value = "Hortaliza"
